$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# New reference-catalog row (row 28), continuing numbering from row 27 (22)
# Shared-string interning order matters: municipalDistrict, then the
# Russian label, then the source link (matches the source file's order).
$ws.Cells.Item(28, 2).Value = 23
$ws.Cells.Item(28, 4).Value = "municipalDistrict"
$ws.Cells.Item(28, 3).Value = "Муниципальные образования"
$ws.Cells.Item(28, 5).Value = "https://rosstat.gov.ru/opendata/7708234640-oktmo"

# Copy the formatting from the row above (row 27) so the new row matches
$ws.Range("B27:G27").Copy()
$ws.Range("B28:G28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Clear F28/G28 text that may have been brought over by the paste (row 27 F/G differ)
$ws.Cells.Item(28, 6).Value = $null
$ws.Cells.Item(28, 7).Value = $null

# Update the saved cursor/selection position
$ws.Range("E20").Select()
